# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# B11 on the "Rules" sheet currently holds the shared string "R40"
# (rule R40's label). The commit changes that cell's content to the
# literal text "1" while leaving the cell's existing style/formatting
# untouched.
#
# A plain `Range.Value = "1"` assignment would be auto-typed as the
# number 1 (Excel parses numeric-looking input), and forcing text via
# NumberFormat="@" changes the cell's style index. To keep the original
# style intact we stage the text value (with a Text number format) in a
# scratch cell, then copy/paste only the *value* onto B11 so its
# formatting is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
